# Update crypto price (D) and volume-change (E) columns per the commit diff.
# Numeric-looking new text values are written with a leading apostrophe
# (quote-prefix) so Excel keeps/stores them as text, matching the source
# data which is textual (e.g. "314.23", "0.100"), not a real number.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "40.963.21"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "2.389.49"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'314.23"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "'88.84"
$ws.Range("E6").Value = "  -4.79%  "
$ws.Range("E7").Value = "  -3.70%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  -4.13%  "
$ws.Range("D10").Value = "'0.0843"
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("D11").Value = "'31.31"
$ws.Range("E11").Value = "  -6.15%  "
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").Value = "2.759.66"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("D14").Value = "'6.62"
$ws.Range("E14").Value = "  -4.44%  "
$ws.Range("D15").Value = "'15.29"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").Value = "2.378.78"
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("E17").Value = "  -3.52%  "
$ws.Range("D18").Value = "40.933.17"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").Value = "0.0₃0917"
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("D20").Value = "'6.19"
$ws.Range("D21").Value = "'69.40"
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("E22").Value = "  -3.97%  "
$ws.Range("D23").Value = "'233.53"
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("E24").Value = "  -3.66%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'1.83"
$ws.Range("E26").Value = "  -6.28%  "
$ws.Range("D27").Value = "'24.08"
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("D29").Value = "'9.40"
$ws.Range("E29").Value = "  -4.09%  "
$ws.Range("D30").Value = "'34.18"
$ws.Range("E30").Value = "  -7.30%  "
$ws.Range("D31").Value = "'153.94"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'5.23"
$ws.Range("E33").Value = "  -5.20%  "
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("E35").Value = "  -4.97%  "
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("D38").Value = "'16.13"
$ws.Range("E38").Value = "  -8.05%  "
$ws.Range("D39").Value = "'0.100"
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("E40").Value = "  -7.44%  "
$ws.Range("D41").Value = "'3.86"
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("E42").Value = "  -7.21%  "
$ws.Range("D43").Value = "1.975.99"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Value = "'0.0272"
$ws.Range("E44").Value = "  -4.60%  "
$ws.Range("D45").Value = "'17.75"
$ws.Range("E45").Value = "  -7.67%  "
$ws.Range("D46").Value = "'9.67"
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("D47").Value = "'2.77"
$ws.Range("E47").Value = "  -7.59%  "
$ws.Range("D48").Value = "2.627.56"
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("D49").Value = "'93.99"
$ws.Range("E49").Value = "  -4.65%  "
$ws.Range("D50").Value = "'72.92"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("D51").Value = "'51.57"
$ws.Range("E51").Value = "  -1.99%  "
